$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ADC resolution rows
$ws.Range("B12").Value = 4096
$ws.Range("B13").Formula = "=C4/B12"

$ws.Range("B15").Value = "ADC Wert"
$ws.Range("C15").Formula = "=C9/B13"
$ws.Range("D15").Formula = "=D9/B13"
$ws.Range("E15").Formula = "=E9/B13"
$ws.Range("F15").Formula = "=F9/B13"
$ws.Range("G15").Formula = "=G9/B13"
$ws.Range("H15").Formula = "=H9/B13"
$ws.Range("I15").Formula = "=I9/B13"
$ws.Range("J15").Formula = "=J9/B13"
$ws.Range("C15:J15").NumberFormat = "0"

# Column widths (best-fit widths computed by Excel when the sheet was authored)
$ws.Columns("C:C").ColumnWidth = 11.3
$ws.Columns("D:J").ColumnWidth = 12.3

# Selection moved to C16 after the edit
$ws.Range("C16").Select()
